$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Match conventions and units with simulink
# E6 was "pBt" -> now "q2p"
$ws.Range("E6").Value = "q2p"

# C5 was "Electrcial_Storage_Units" -> now "Thermal_Storage_Units"
$ws.Range("C5").Value = "Thermal_Storage_Units"

# Update the selected/active cell to C5
$ws.Range("C5").Select()
